$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block for "Форма "Текущее состояние и местоположение техники"" ---
# Row 9: red-fill marker cell (B9) + task text (C9) - continuation of the
# existing task list style used by rows 3-8.
$ws.Range("B3").Copy($ws.Range("B9"))

# Row 12 text is written before row 9's text so the shared-strings table
# ends up in the same order as the source workbook (index 7 = "Форма...",
# index 8 = "Добавить...").
$ws.Range("C12").Value = 'Форма "Текущее состояние и местоположение техники"'

$ws.Range("C9").Value = "Добавить строку управления таблицами"

# Row 11: new date header, re-using the existing date cell's number format.
$ws.Range("B2").Copy($ws.Range("B11"))
$ws.Range("B11").Value = 45053

# Rows 13-18: remaining task descriptions.
$ws.Range("C13").Value = 'Создать справочник "Место расположения техники"'
$ws.Range("C14").Value = 'Создать справочник "Ответсвенные лица"'
$ws.Range("C15").Value = 'Создать документ "Потребность"'
$ws.Range("C16").Value = 'Создать документ "Закупка"'
$ws.Range("C17").Value = 'Создать документ "Списание"'
$ws.Range("C18").Value = 'Создать регистр "Движения техники" Форма "Текущее состояние и местоположение техники"'

# Update the view/selection to match the post-edit state (last selected
# cell after entering the new rows).
$ws.Range("C19").Select()
